$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-7: date serial 45207 -> 45208 (2023-10-08 -> 2023-10-09)
$ws.Range("C2:C7").Value = 45208
